$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New word entries to append (word, definition, example1, example2, rowHeight)
$entries = @(
    @("idiot", "a stupid person or someone who has done something stupid", "I smile like an idiot when i'm talking to you.", "It was all your fault, you idiot.", 45),
    @("immunity", "the state or right of being protected from particular laws or from unpleasant things", "the caccine provides longer immunity against flu.", "They were granted immunity from prosecution.", 60),
    @("narrow", "measuring only a small distance from one side to the other, especially in relation to the length", "a long narrow road", "The stairs were very narrow.", 60),
    @("lane", "a narrow road in the countryside", "the police opened a lane through the crowd and let us pass.", "a quiet country lane", 60),
    @("liable", "legally responsible for the cost of something", "he claimed , he was not liable for his wife's debts.", "You’re more liable to injury when you don’t get regular exercise.", 45),
    @("debt", "a sum of money that a person or organization owes", "She had debts of over £100,000.", "The band will be in debt to the record company for years.", 45),
    @("obliterate", "to destroy something completely so that nothing remains", "the bomb nearly obliterate the city.", "Hiroshima was nearly obliterated by the atomic bomb.", 45),
    @("promenade", "a wide road next to the beach, where people can walk for pleasure", "we look a promenade along the canal after sunday dinner", "This pier was not a promenade for me.", 60)
)

$startRow = 92
for ($i = 0; $i -lt $entries.Count; $i++) {
    $row = $startRow + $i
    $entry = $entries[$i]
    $ws.Range("A$row").Value = $entry[0]
    $ws.Range("B$row").Value = $entry[1]
    $ws.Range("C$row").Value = $entry[2]
    $ws.Range("D$row").Value = $entry[3]
    $ws.Rows.Item($row).RowHeight = $entry[4]
}

# row 96 column A keeps no explicit style (matches source quirk)
$ws.Range("A96").Style = "Normal"

# Row 91's height was adjusted as part of this edit too
$ws.Rows.Item(91).RowHeight = 60

# Update sheet view (scroll position / selection) to match post-edit state
$ws.Application.ActiveWindow.ScrollRow = 94
$sel = $ws.Range("C102")
$sel.Select()
